$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values that look like pure numbers but must stay as text (matching the
# original inline-string cell type) are written via NumberFormat "@" with the
# original Style restored afterwards, so the cell keeps its original appearance.

$ws.Range('D2').Value = '37.707.04'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').Value = '2.078.46'
$ws.Range('E3').Value = '  +4.31%  '
$ws.Range('E4').Value = '  -0.06%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.08'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  -1.92%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.617'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  +2.21%  '
$origStyle = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.45'
$ws.Range('D7').Style = $origStyle
$ws.Range('E7').Value = '  +7.75%  '
$ws.Range('E8').Value = '  +0.02%  '
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.385'
$ws.Range('D9').Style = $origStyle
$ws.Range('E9').Value = '  +3.14%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.16'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  +1.82%  '
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0765'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  +1.60%  '
$ws.Range('E12').Value = '  +3.98%  '
$ws.Range('D13').Value = '2.385.81'
$ws.Range('E13').Value = '  +4.27%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.47'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  +2.18%  '
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.05'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  +0.33%  '
$origStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.782'
$ws.Range('D16').Style = $origStyle
$ws.Range('E16').Value = '  +3.29%  '
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.24'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  +3.60%  '
$ws.Range('D18').Value = '2.054.08'
$ws.Range('E18').Value = '  +2.84%  '
$ws.Range('D19').Value = '37.883.55'
$ws.Range('E19').Value = '  +2.71%  '
$ws.Range('E20').Value = '  +21.15%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '68.87'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  +0.56%  '
$ws.Range('D22').Value = '0.0₃0815'
$ws.Range('E22').Value = '  +0.80%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '225.21'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  -1.23%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  -0.06%  '
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.46'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  +1.27%  '
$ws.Range('E26').Value = '  +2.81%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.75'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  +0.63%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.87'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  +2.40%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.131'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  +3.35%  '
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.38'
$ws.Range('D30').Style = $origStyle
$ws.Range('E31').Value = '  +6.57%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.118'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +0.67%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0632'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  +3.53%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.49'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  +1.40%  '
$ws.Range('E35').Value = '  +13.30%  '
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.48'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  +5.84%  '
$ws.Range('E37').Value = '  -0.19%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.96'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +12.13%  '
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.35'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  +2.76%  '
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('E41').Value = '  -3.01%  '
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0964'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  +7.30%  '
$ws.Range('D43').Value = '1.485.35'
$ws.Range('E43').Value = '  +4.01%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '95.64'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  +8.57%  '
$ws.Range('E45').Value = '  +4.33%  '
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.48'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  +9.51%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.25'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  +26.42%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.13'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  +0.88%  '
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.36'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  +9.48%  '
$ws.Range('E50').Value = '  +2.28%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.93'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  +1.69%  '
